$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header row (row 1): "Edad" -> "Fecha Matrícula", "Modalidad" -> "Correo"
# ---------------------------------------------------------------------------
$ws.Cells.Item(1, 3).Value = "Fecha Matrícula"
$ws.Cells.Item(1, 4).Value = "Correo"

# ---------------------------------------------------------------------------
# Existing data row (row 2): Edad -> Fecha Matricula, Modalidad -> Correo
# A leading apostrophe forces the value to be kept as literal text instead of
# being auto-converted to a date/number.
# ---------------------------------------------------------------------------
$ws.Cells.Item(2, 3).Formula = "'2025-01-04"
$ws.Cells.Item(2, 4).Value = "paola@gmail.com"

# Re-apply the same look (alignment / font / border) that the rest of row 2
# already has, so the edited cell still matches its neighbours.
$ws.Range("C2").HorizontalAlignment = $ws.Range("A2").HorizontalAlignment
$ws.Range("C2").VerticalAlignment = $ws.Range("A2").VerticalAlignment
$ws.Range("C2").Font.Bold = $ws.Range("A2").Font.Bold
$ws.Range("C2").Borders.LineStyle = $ws.Range("A2").Borders.LineStyle

# ---------------------------------------------------------------------------
# Remove column E (old "Peso" column) entirely - data + column definition
# ---------------------------------------------------------------------------
$ws.Columns.Item(5).Delete()

# ---------------------------------------------------------------------------
# New data row (row 3)
# ---------------------------------------------------------------------------
$ws.Cells.Item(3, 1).Formula = "'1716776412"
$ws.Cells.Item(3, 2).Value = "hfgfghfgh fghfghfgh"
$ws.Cells.Item(3, 3).Formula = "'2025-02-23"
$ws.Cells.Item(3, 4).Value = "darwi@ngmail.com"

# Give row 3 the same look & feel (style) as row 2
$ws.Range("A3:D3").HorizontalAlignment = $ws.Range("A2:D2").HorizontalAlignment
$ws.Range("A3:D3").VerticalAlignment = $ws.Range("A2:D2").VerticalAlignment
$ws.Range("A3:D3").Font.Bold = $ws.Range("A2:D2").Font.Bold
$ws.Range("A3:D3").Borders.LineStyle = $ws.Range("A2:D2").Borders.LineStyle

# ---------------------------------------------------------------------------
# Column widths for C and D (the values below account for Excel's internal
# char-width/pixel rounding so the saved width lands exactly on 17 / 18)
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 16.15
$ws.Columns.Item(4).ColumnWidth = 17.15
